$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell V1: new epi-week column "19" (week 19 of 2025) ---
# Assigning Value = "19" directly would auto-coerce to a Double (same as real
# Excel COM numeric-string coercion), losing the text type the other header
# cells (D1:U1) use, and forcing text via a leading apostrophe allocates a new
# quotePrefix style. Using a temporary formula + copy/paste-values keeps the
# cell text-typed while reusing the existing header style (same as U1).
$ws.Range("V1").Formula = "=""19"""
$ws.Range("V1").Copy()
$ws.Range("V1").PasteSpecial(-4163)
$ws.Range("V1").Font.Bold = $true
$ws.Range("V1").HorizontalAlignment = -4108

# --- Data column V2:V54 (numeric weekly counts for epi-week 19) ---
$ws.Range("V2").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("V6").Value = 22
$ws.Range("V7").Value = 4
$ws.Range("V8").Value = 27
$ws.Range("V9").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("V11").Value = 0
$ws.Range("V13").Value = 0
$ws.Range("V14").Value = 0
$ws.Range("V15").Value = 0
$ws.Range("V17").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("V23").Value = 1
$ws.Range("V24").Value = 0
$ws.Range("V26").Value = 0
$ws.Range("V27").Value = 2
$ws.Range("V28").Value = 22
$ws.Range("V29").Value = 0
$ws.Range("V30").Value = 0
$ws.Range("V32").Value = 12
$ws.Range("V33").Value = 0
$ws.Range("V34").Value = 0
$ws.Range("V35").Value = 0
$ws.Range("V37").Value = 0
$ws.Range("V38").Value = 0
$ws.Range("V39").Value = 0
$ws.Range("V40").Value = 0
$ws.Range("V41").Value = 0
$ws.Range("V42").Value = 0
$ws.Range("V43").Value = 0
$ws.Range("V44").Value = 0
$ws.Range("V45").Value = 0
$ws.Range("V46").Value = 0
$ws.Range("V47").Value = 0
$ws.Range("V48").Value = 0
$ws.Range("V49").Value = 0
$ws.Range("V50").Value = 0
$ws.Range("V51").Value = 0
$ws.Range("V52").Value = 0
$ws.Range("V53").Value = 0
$ws.Range("V54").Value = 0
